$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell $ws "B2" "Bitcoin"
Set-TextCell $ws "C2" "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell $ws "D2" "29.410.81"
Set-TextCell $ws "E2" "  -1.58%  "
Set-TextCell $ws "B3" "Ethereum"
Set-TextCell $ws "C3" "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell $ws "D3" "1.849.53"
Set-TextCell $ws "E3" "  -0.50%  "
Set-TextCell $ws "B4" "TetherUSD"
Set-TextCell $ws "C4" "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell $ws "D4" "0.9995"
Set-TextCell $ws "E4" "  -0.03%  "
Set-TextCell $ws "B5" "BNB"
Set-TextCell $ws "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell $ws "D5" "243.00"
Set-TextCell $ws "E5" "  -1.36%  "
Set-TextCell $ws "B6" "XRP"
Set-TextCell $ws "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell $ws "D6" "0.6576"
Set-TextCell $ws "E6" "  +3.18%  "
Set-TextCell $ws "B7" "USDC"
Set-TextCell $ws "C7" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell $ws "D7" "1.000"
Set-TextCell $ws "E7" "  +0.01%  "
Set-TextCell $ws "B8" "OKB"
Set-TextCell $ws "C8" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D8" "48.04"
Set-TextCell $ws "E8" "  +2.83%  "
Set-TextCell $ws "B9" "Cardano"
Set-TextCell $ws "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell $ws "D9" "0.2987"
Set-TextCell $ws "E9" "  -0.37%  "
Set-TextCell $ws "B10" "Dogecoin"
Set-TextCell $ws "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell $ws "D10" "0.07473"
Set-TextCell $ws "E10" "  +0.04%  "
Set-TextCell $ws "B11" "Solana"
Set-TextCell $ws "C11" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell $ws "D11" "24.31"
Set-TextCell $ws "E11" "  -0.84%  "
Set-TextCell $ws "B12" "TRON"
Set-TextCell $ws "C12" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D12" "0.07629"
Set-TextCell $ws "E12" "  -0.62%  "
Set-TextCell $ws "B13" "WrappedEther"
Set-TextCell $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D13" "1.849.19"
Set-TextCell $ws "E13" "  -0.58%  "
Set-TextCell $ws "B14" "Polkadot"
Set-TextCell $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D14" "5.015"
Set-TextCell $ws "E14" "  -0.66%  "
Set-TextCell $ws "B15" "Polygon"
Set-TextCell $ws "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell $ws "D15" "0.6840"
Set-TextCell $ws "E15" "  -1.06%  "
Set-TextCell $ws "B16" "Litecoin"
Set-TextCell $ws "C16" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws "D16" "83.59"
Set-TextCell $ws "E16" "  -0.85%  "
Set-TextCell $ws "B17" "ShibaInu"
Set-TextCell $ws "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D17" "0.000009480"
Set-TextCell $ws "E17" "  +1.53%  "
Set-TextCell $ws "B18" "Uniswap"
Set-TextCell $ws "C18" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws "D18" "6.141"
Set-TextCell $ws "E18" "  +1.11%  "
Set-TextCell $ws "B19" "WrappedBTC"
Set-TextCell $ws "C19" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws "D19" "29.451.01"
Set-TextCell $ws "E19" "  -1.34%  "
Set-TextCell $ws "B20" "WrappedliquidstakedEther2.0"
Set-TextCell $ws "C20" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell $ws "D20" "2.067.33"
Set-TextCell $ws "E20" "  -2.39%  "
Set-TextCell $ws "B21" "BitcoinCash"
Set-TextCell $ws "C21" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D21" "237.22"
Set-TextCell $ws "E21" "  -0.53%  "
Set-TextCell $ws "B22" "Avalanche"
Set-TextCell $ws "C22" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws "D22" "12.56"
Set-TextCell $ws "E22" "  -0.82%  "
Set-TextCell $ws "B23" "Dai"
Set-TextCell $ws "C23" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D23" "1.000"
Set-TextCell $ws "E23" "  -0.01%  "
Set-TextCell $ws "B24" "Chainlink"
Set-TextCell $ws "C24" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D24" "7.655"
Set-TextCell $ws "E24" "  +3.95%  "
Set-TextCell $ws "B25" "BinanceUSD"
Set-TextCell $ws "C25" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell $ws "D25" "1.000"
Set-TextCell $ws "E25" "  -0.07%  "
Set-TextCell $ws "B26" "Stellar"
Set-TextCell $ws "C26" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D26" "0.1425"
Set-TextCell $ws "E26" "  +0.58%  "
Set-TextCell $ws "B27" "Monero"
Set-TextCell $ws "C27" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D27" "156.81"
Set-TextCell $ws "E27" "  -1.45%  "
Set-TextCell $ws "B28" "Cosmos"
Set-TextCell $ws "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D28" "8.478"
Set-TextCell $ws "E28" "  -1.16%  "
Set-TextCell $ws "B29" "EthereumClassic"
Set-TextCell $ws "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D29" "17.78"
Set-TextCell $ws "E29" "  -1.02%  "
Set-TextCell $ws "B30" "Hedera"
Set-TextCell $ws "C30" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D30" "0.06034"
Set-TextCell $ws "E30" "  -0.45%  "
Set-TextCell $ws "B31" "PancakeSwap"
Set-TextCell $ws "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D31" "1.488"
Set-TextCell $ws "E31" "  -1.09%  "
Set-TextCell $ws "B32" "Toncoin"
Set-TextCell $ws "C32" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D32" "1.251"
Set-TextCell $ws "E32" "  -1.89%  "
Set-TextCell $ws "B33" "Filecoin"
Set-TextCell $ws "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D33" "4.135"
Set-TextCell $ws "E33" "  -0.16%  "
Set-TextCell $ws "B34" "InternetComputer(DFINITY)"
Set-TextCell $ws "C34" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D34" "4.072"
Set-TextCell $ws "E34" "  -1.42%  "
Set-TextCell $ws "B35" "LidoDAOToken"
Set-TextCell $ws "C35" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D35" "1.856"
Set-TextCell $ws "E35" "  -1.80%  "
Set-TextCell $ws "B36" "ARBITRUM"
Set-TextCell $ws "C36" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D36" "1.175"
Set-TextCell $ws "E36" "  +0.77%  "
Set-TextCell $ws "B37" "ImmutableX"
Set-TextCell $ws "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D37" "0.7207"
Set-TextCell $ws "E37" "  -1.09%  "
Set-TextCell $ws "B38" "HuobiToken"
Set-TextCell $ws "C38" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws "D38" "2.597"
Set-TextCell $ws "E38" "  -0.52%  "
Set-TextCell $ws "B39" "MXToken"
Set-TextCell $ws "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D39" "2.792"
Set-TextCell $ws "E39" "  -2.28%  "
Set-TextCell $ws "B40" "VeChain"
Set-TextCell $ws "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D40" "0.01781"
Set-TextCell $ws "E40" "  -0.86%  "
Set-TextCell $ws "B41" "Maker"
Set-TextCell $ws "C41" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws "D41" "1.197.23"
Set-TextCell $ws "E41" "  -2.17%  "
Set-TextCell $ws "B42" "TrustWalletToken"
Set-TextCell $ws "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D42" "0.9117"
Set-TextCell $ws "E42" "  -2.61%  "
Set-TextCell $ws "B43" "FraxShare"
Set-TextCell $ws "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D43" "6.225"
Set-TextCell $ws "E43" "  -0.74%  "
Set-TextCell $ws "B44" "PaxDollar"
Set-TextCell $ws "C44" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws "D44" "1.000"
Set-TextCell $ws "E44" "  -0.17%  "
Set-TextCell $ws "B45" "RocketPoolETH"
Set-TextCell $ws "C45" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws "D45" "2.002.25"
Set-TextCell $ws "E45" "  -1.35%  "
Set-TextCell $ws "B46" "Quant"
Set-TextCell $ws "C46" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws "D46" "101.92"
Set-TextCell $ws "E46" "  -0.35%  "
Set-TextCell $ws "B47" "Aave"
Set-TextCell $ws "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws "D47" "66.21"
Set-TextCell $ws "E47" "  -0.02%  "
Set-TextCell $ws "B48" "Aptos"
Set-TextCell $ws "C48" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws "D48" "7.434"
Set-TextCell $ws "E48" "  +10.76%  "
Set-TextCell $ws "B49" "BabyDogeCoin"
Set-TextCell $ws "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws "D49" "0.00000000122"
Set-TextCell $ws "E49" "  -0.27%  "
Set-TextCell $ws "B50" "TheSandbox"
Set-TextCell $ws "C50" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws "D50" "0.4048"
Set-TextCell $ws "E50" "  -1.03%  "
Set-TextCell $ws "B51" "EnergySwap"
Set-TextCell $ws "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D51" "9.083"
Set-TextCell $ws "E51" "  -2.63%  "
